$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New list of node labels (grew from i1..i4 to i1..i9)
$labels = @("i1","i2","i3","i4","i5","i6","i7","i8","i9")

# Write header row (B1:J1)
for ($j = 0; $j -lt $labels.Length; $j++) {
    $ws.Cells.Item(1, 2 + $j).Value = $labels[$j]
}

# Write row labels (A2:A10)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $labels[$i]
}

# Symmetric matrix of values for the 9x9 block (rows 2..10, cols B..J)
$matrix = @(
    @(0, 0.00178640640744647, 0.001783378767856305, 0.001183898155703645, 0, 0, 0, 0, 0),
    @(0.00178640640744647, 0, 0.001436737742693641, 5.786744453957676, 0, 0, 0, 0, 0),
    @(0.001783378767856305, 0.001436737742693641, 0, 2.257054535163327, 0, 0, 0, 0, 0),
    @(0.001183898155703645, 5.786744453957676, 2.257054535163327, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0)
)

for ($i = 0; $i -lt 9; $i++) {
    $row = $matrix[$i]
    for ($j = 0; $j -lt 9; $j++) {
        $ws.Cells.Item(2 + $i, 2 + $j).Value = $row[$j]
    }
}
